# Translate the Sheet1 header row (row 1) from the original Spanish
# labels to the new English machine-readable codes. Overwriting these
# cells causes the old Spanish header strings to fall out of the shared
# strings table (they are no longer referenced anywhere else in the
# sheet) while the new English strings are appended - this naturally
# reproduces the shared-strings renumbering seen in the diff, with no
# other cell value in the sheet actually changing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "report_number",
    "so",
    "assignee",
    "calibrator",
    "supervisor",
    "approver",
    "classification",
    "lab_received_date",
    "scheduled_delivery_date",
    "cleaning_date",
    "calibration_date",
    "delivery_date",
    "delivery_time",
    "process_status",
    "assigned_time",
    "service_location",
    "substitution_reason",
    "lab_observations",
    "priority",
    "piece_count"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 18 holds a multi-line "NINGUNA" note in column R; give it the
# taller row height Excel computed for the wrapped text.
$ws.Rows.Item(18).RowHeight = 28.8

# Recreate the hidden _FilterDatabase defined name (sheet-scoped) that
# Excel writes once a filter has been defined over column B.
$ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$B`$1:`$B`$37")
$fdb = $ws.Names.Item($ws.Names.Count)
$fdb.Visible = $false

# Restore the last-saved selection.
$ws.Range("R23").Select()
